$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell F1 with the same style as the other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$timestamps = @(
    "2021-10-05 13:38:35.083788",
    "2021-10-05 13:38:35.083799",
    "2021-10-05 13:38:35.083803",
    "2021-10-05 13:38:35.083806",
    "2021-10-05 13:38:35.083809",
    "2021-10-05 13:38:35.083811",
    "2021-10-05 13:38:35.083814",
    "2021-10-05 13:38:35.083817",
    "2021-10-05 13:38:35.083819",
    "2021-10-05 13:38:35.083822",
    "2021-10-05 13:38:35.083824",
    "2021-10-05 13:38:35.083827",
    "2021-10-05 13:38:35.083830",
    "2021-10-05 13:38:35.083832",
    "2021-10-05 13:38:35.083835",
    "2021-10-05 13:38:35.083837",
    "2021-10-05 13:38:35.083840",
    "2021-10-05 13:38:35.083843",
    "2021-10-05 13:38:35.083845",
    "2021-10-05 13:38:35.083848",
    "2021-10-05 13:38:35.083850",
    "2021-10-05 13:38:35.083853",
    "2021-10-05 13:38:35.083855",
    "2021-10-05 13:38:35.083858",
    "2021-10-05 13:38:35.083861",
    "2021-10-05 13:38:35.083863",
    "2021-10-05 13:38:35.083866",
    "2021-10-05 13:38:35.083869",
    "2021-10-05 13:38:35.083871",
    "2021-10-05 13:38:35.083874",
    "2021-10-05 13:38:35.083876",
    "2021-10-05 13:38:35.083879",
    "2021-10-05 13:38:35.083881",
    "2021-10-05 13:38:35.083884",
    "2021-10-05 13:38:35.083887",
    "2021-10-05 13:38:35.083889",
    "2021-10-05 13:38:35.083892",
    "2021-10-05 13:38:35.083895",
    "2021-10-05 13:38:35.083897",
    "2021-10-05 13:38:35.083900",
    "2021-10-05 13:38:35.083903",
    "2021-10-05 13:38:35.083906",
    "2021-10-05 13:38:35.083908",
    "2021-10-05 13:38:35.083911",
    "2021-10-05 13:38:35.083913",
    "2021-10-05 13:38:35.083916",
    "2021-10-05 13:38:35.083919",
    "2021-10-05 13:38:35.083921",
    "2021-10-05 13:38:35.083924",
    "2021-10-05 13:38:35.083926",
    "2021-10-05 13:38:35.083929",
    "2021-10-05 13:38:35.083931",
    "2021-10-05 13:38:35.083934",
    "2021-10-05 13:38:35.083937",
    "2021-10-05 13:38:35.083939",
    "2021-10-05 13:38:35.083942",
    "2021-10-05 13:38:35.083944",
    "2021-10-05 13:38:35.083947",
    "2021-10-05 13:38:35.083949",
    "2021-10-05 13:38:35.083952",
    "2021-10-05 13:38:35.083954",
    "2021-10-05 13:38:35.083957",
    "2021-10-05 13:38:35.083960",
    "2021-10-05 13:38:35.083962",
    "2021-10-05 13:38:35.083966",
    "2021-10-05 13:38:35.083969",
    "2021-10-05 13:38:35.083972",
    "2021-10-05 13:38:35.083975",
    "2021-10-05 13:38:35.083977",
    "2021-10-05 13:38:35.083980",
    "2021-10-05 13:38:35.083982",
    "2021-10-05 13:38:35.083985",
    "2021-10-05 13:38:35.083987",
    "2021-10-05 13:38:35.083990",
    "2021-10-05 13:38:35.083993",
    "2021-10-05 13:38:35.083995",
    "2021-10-05 13:38:35.084000",
    "2021-10-05 13:38:35.084003",
    "2021-10-05 13:38:35.084006",
    "2021-10-05 13:38:35.084009",
    "2021-10-05 13:38:35.084011",
    "2021-10-05 13:38:35.084014",
    "2021-10-05 13:38:35.084017",
    "2021-10-05 13:38:35.084020",
    "2021-10-05 13:38:35.084023",
    "2021-10-05 13:38:35.084025",
    "2021-10-05 13:38:35.084028",
    "2021-10-05 13:38:35.084031",
    "2021-10-05 13:38:35.084033",
    "2021-10-05 13:38:35.084036",
    "2021-10-05 13:38:35.084039",
    "2021-10-05 13:38:35.084041",
    "2021-10-05 13:38:35.084046",
    "2021-10-05 13:38:35.084049",
    "2021-10-05 13:38:35.084051",
    "2021-10-05 13:38:35.084054",
    "2021-10-05 13:38:35.084057",
    "2021-10-05 13:38:35.084059",
    "2021-10-05 13:38:35.084062",
    "2021-10-05 13:38:35.084065",
    "2021-10-05 13:38:35.084067",
    "2021-10-05 13:38:35.084070",
    "2021-10-05 13:38:35.084072",
    "2021-10-05 13:38:35.084075",
    "2021-10-05 13:38:35.084078"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
